$wb = $excel.ActiveWorkbook

# --- "About" sheet: add a Notes section at A48:A50 ---
$about = $wb.Worksheets.Item("About")

$about.Range("A48").Value = "Notes"
$about.Range("A48").Font.Bold = $true

$about.Range("A49").Value = 'This variable is also known as the "Fuel Economy Rebound Effect" or "Fuel Cost Rebound Effect." It is the change'
$about.Range("A50").Value = "in VMT as a fraction of the change in fuel cost. E.g. for a 1% increase in fuel cost per mile, VMT changes by -0.1%."

# --- "EoDfVUwFC" sheet: clarify the Elasticity header ---
$eod = $wb.Worksheets.Item("EoDfVUwFC")
$eod.Activate()

$eod.Range("B1").Value = "Elasticity (dimensionless)"
$eod.Range("B1").WrapText = $true
$eod.Rows.Item(1).RowHeight = 30

$eod.Range("B1").Select()

# Return focus to the "About" sheet, scrolled down to the new notes section
$about.Activate()
$about.Range("A48").Select()
